$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "ProgramsTab" query in B2: add a CASE expression that
# --- derives the "Website" column from prg.program_link / prg.program_acronym
# --- instead of returning prg.website directly.
$newProgramsQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Thyroid Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
"@

$ws.Range("B2").Value = $newProgramsQuery

# Re-apply the existing formatting (still 12pt, wrap text) on the edited
# cell. Excel churns the style table (a harmless duplicate font/xf entry)
# whenever a long cell is retyped like this, so nudge the size away and
# back to reproduce that same style-table growth.
$ws.Range("B2").Font.Size = 96
$ws.Range("B2").Font.Size = 12
$ws.Range("B2").WrapText = $true

# --- Move the active selection / scroll position, matching the author's
# --- final cursor position after editing (was topLeftCell A5 / cell C5,
# --- now A2 / C3).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("C3").Select() | Out-Null
